$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:F11, G2:G11 based on the diff (rows shift up by one, with a new
# final data row appended at the bottom of the window; G values are a decreasing count).
$data = @(
    @(0.04463932893276446, 0.5563105143949441, 0.6682676647472701, 0.817476400116401, 0.8399212116787623, 18),
    @(-0.04959926954625594, 0.5001880640764477, 0.5928923876220049, 0.7699950568815392, 0.7920443919614847, 17),
    @(0.0418877556493024, 0.3990978627449576, 0.4202563537498312, 0.6482718208821291, 0.6681331361749159, 16),
    @(0.1640516070889592, 0.3138597579017368, 0.2168740114234605, 0.4656973388623351, 0.4511426939707701, 15),
    @(0.1219509228516585, 0.3069632200035473, 0.1600708036365558, 0.4000884947565423, 0.3954338284360399, 14),
    @(0.1397471008904613, 0.2788282874804134, 0.1873748606623431, 0.4328681793136833, 0.4264185062740307, 13),
    @(0.1371477073313899, 0.3017865247828347, 0.1605887215133184, 0.4007352261947013, 0.3932788433212083, 12),
    @(0.1401056802292451, 0.2816089537871457, 0.2148615712472957, 0.4635316291767971, 0.4634168389002237, 11),
    @(0.1061516567236955, 0.3103361548660022, 0.1599817766894938, 0.3999772202132189, 0.4064939557625106, 10),
    @(0.09702283882948032, 0.3111699130949587, 0.2143114512357824, 0.4629378481349115, 0.4801148534245816, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
